$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 210, pushing rows 210-255 down to 212-257
$ws.Rows.Item(210).Resize(2).Insert()

# New row 210 — "Especial" quality, 10-unit box
$ws.Cells.Item(210, 1).Value = 4
$ws.Cells.Item(210, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(210, 3).Value = "Los Lagos"
$ws.Cells.Item(210, 4).Value = 44722
$ws.Cells.Item(210, 5).Value = 10
$ws.Cells.Item(210, 6).Value = "Fruta"
$ws.Cells.Item(210, 7).Value = 100108
$ws.Cells.Item(210, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(210, 9).Value = 100108005
$ws.Cells.Item(210, 10).Value = "Piña"
$ws.Cells.Item(210, 11).Value = "Caramelo"
$ws.Cells.Item(210, 12).Value = "Especial"
$ws.Cells.Item(210, 13).Value = 200
$ws.Cells.Item(210, 14).Value = 18000
$ws.Cells.Item(210, 15).Value = 19000
$ws.Cells.Item(210, 16).Value = 18500
$ws.Cells.Item(210, 17).Value = "`$/caja 10 unidades"
$ws.Cells.Item(210, 18).Value = "Ecuador"
$ws.Cells.Item(210, 19).Value = 1850
$ws.Cells.Item(210, 20).Value = 10

# New row 211 — "Tercera" quality, 16-unit box
$ws.Cells.Item(211, 1).Value = 4
$ws.Cells.Item(211, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(211, 3).Value = "Los Lagos"
$ws.Cells.Item(211, 4).Value = 44722
$ws.Cells.Item(211, 5).Value = 10
$ws.Cells.Item(211, 6).Value = "Fruta"
$ws.Cells.Item(211, 7).Value = 100108
$ws.Cells.Item(211, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(211, 9).Value = 100108005
$ws.Cells.Item(211, 10).Value = "Piña"
$ws.Cells.Item(211, 11).Value = "Caramelo"
$ws.Cells.Item(211, 12).Value = "Tercera"
$ws.Cells.Item(211, 13).Value = 200
$ws.Cells.Item(211, 14).Value = 20000
$ws.Cells.Item(211, 15).Value = 21000
$ws.Cells.Item(211, 16).Value = 20500
$ws.Cells.Item(211, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(211, 18).Value = "Ecuador"
$ws.Cells.Item(211, 19).Value = 1281
$ws.Cells.Item(211, 20).Value = 16
